$d = $word.ActiveDocument

# Step 1: replace text of the first bullet (relay button) with the new merged text
$find = $d.Content.Find
$find.Execute("Berhasil membangun project sederhana simulasi relay button", $true, $false, $false, $false, $false, $true, 1, $false, "Berhasil mengakses API melalui simulasi Wokwi menggunakan Ngrok, dan Laravel", 2)

# Step 2: delete the entire second bullet paragraph (ultrasonic) including its paragraph mark
$p = $d.Paragraphs.Item(31)
$p.Range.Delete()

for ($i = 28; $i -le 32; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output "$i => [$($pp.Range.Text)]"
}
